$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.887.04"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.91%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.520"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0900"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "1.868.20"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "1.629.73"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.587"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.64%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.907.85"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").Value = "1.424.25"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.73%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0504"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.835"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.775.84"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.57%  "
